$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 0.004000186920166016

$ws.Range("C3").Value = 640.41
$ws.Range("D3").Value = 493.05
$ws.Range("E3").Value = 147.36
$ws.Range("F3").Value = 1133.46
$ws.Range("G3").Value = 566.73
$ws.Range("H3").Value = 1.642164468765259

$ws.Range("H4").Value = 0.003081560134887695

$ws.Range("C5").Value = 515.24
$ws.Range("D5").Value = 512.42
$ws.Range("E5").Value = 2.82
$ws.Range("F5").Value = 1027.66
$ws.Range("G5").Value = 513.83
$ws.Range("H5").Value = 1.55881142616272

$ws.Range("H6").Value = 0.00599980354309082

$ws.Range("C7").Value = 555.13
$ws.Range("D7").Value = 549.67
$ws.Range("E7").Value = 5.47
$ws.Range("F7").Value = 1657.6
$ws.Range("G7").Value = 552.53
$ws.Range("H7").Value = 3.508722305297852

$ws.Range("H8").Value = 0.01004457473754883

$ws.Range("C9").Value = 700.83
$ws.Range("D9").Value = 454.02
$ws.Range("E9").Value = 246.81
$ws.Range("F9").Value = 2322.92
$ws.Range("G9").Value = 580.73
$ws.Range("H9").Value = 6.381536722183228

$ws.Range("H10").Value = 0.01105976104736328

$ws.Range("C11").Value = 540.55
$ws.Range("D11").Value = 530.55
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 2140.9
$ws.Range("G11").Value = 535.22
$ws.Range("H11").Value = 5.761125564575195
